$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

# --- 1. Collapse the three "positioning/sizing/scale" bullets into one ---
$content = $s.Shapes.Item(2)
if ($content.Name -ne "Content") {
    $content = $s.Shapes.Item(3)
}
$content.TextFrame.TextRange.Text = "PNG, JPG, GIF, BMP, TIFF formats`rCustom positioning and sizing`rAspect ratio preservation"

# --- 2. Add two image-placeholder shapes ---
# The slide currently owns shape ids 1,2,3. PowerPoint COM hands out the next
# free id per AddShape call and never reuses one, so burn through ids 4..19
# with throw-away shapes, delete them, and the next two real AddShape calls
# land on ids 20 and 21 (matching the target OOXML `cNvPr id="20"/"21"`).
while ($s.Shapes.Count -lt 18) {
    $tmp = $s.Shapes.AddShape(1, 0, 0, 1, 1)
}
while ($s.Shapes.Count -gt 2) {
    $s.Shapes.Item(3).Delete()
}

function Add-ImagePlaceholder($slide, $name, $offX, $offY, $extCx, $extCy, $label) {
    $leftPt = $offX / 12700.0
    $topPt = $offY / 12700.0
    $widthPt = $extCx / 12700.0
    $heightPt = $extCy / 12700.0

    $shp = $slide.Shapes.AddShape(1, $leftPt, $topPt, $widthPt, $heightPt)
    $shp.Name = $name

    $shp.Fill.ForeColor.RGB = 14737632
    $shp.Line.ForeColor.RGB = 8421504
    $shp.Line.Weight = 1

    $shp.TextFrame.WordWrap = -1
    $shp.TextFrame.VerticalAnchor = 3

    $tr = $shp.TextFrame.TextRange
    $tr.Text = $label
    $tr.ParagraphFormat.Alignment = 2
    $tr.Font.Size = 14

    return $shp
}

$logo = Add-ImagePlaceholder $s "Image Placeholder: logo.png" 100000 100000 2000000 1000000 "📷 logo.png"
$photo = Add-ImagePlaceholder $s "Image Placeholder: photo.jpg" 300000 200000 4000000 3000000 "📷 photo.jpg"

Write-Host "logo id:" $logo.Id
Write-Host "photo id:" $photo.Id
